$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) cells that hold numeric-looking text keep their
# original text representation (e.g. trailing zeros, thousand-separator dots)
# by forcing a Text number format before writing the new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.919.48"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.703.21"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.67"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4085"
$ws.Range("E7").Value = "  +3.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4064"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.003"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.76"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.468"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08828"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.81"
$ws.Range("E13").Value = "  +5.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.497"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.051"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001352"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.635.39"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.55"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07174"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.96"
$ws.Range("E20").Value = "  +5.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.232"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.56"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.929.29"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.884"
$ws.Range("E26").Value = "  -4.95%  "
$ws.Range("B27").Value = "HuobiToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.701"
$ws.Range("E27").Value = "  +28.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.04"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.39"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.77"
$ws.Range("E30").Value = "  +4.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.244"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.280"
$ws.Range("E32").Value = "  +14.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.926.03"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08760"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03209"
$ws.Range("E35").Value = "  +11.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.327"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.020"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8516"
$ws.Range("E39").Value = "  +8.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.92"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09443"
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.06"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.68"
$ws.Range("E44").Value = "  +6.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.728"
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7447"
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.238"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.404"
$ws.Range("E48").Value = "  +6.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.79"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08354"
$ws.Range("E51").Value = "  +4.42%  "
